$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Random Forest (row 4) and Gradient Boosting Machine (row 5) results ---
$ws.Cells.Item(4, 2).Value = 0.155
$ws.Cells.Item(4, 3).Value = 0.129
$ws.Cells.Item(4, 4).Value = 0.989

$ws.Cells.Item(5, 2).Value = 0.211
$ws.Cells.Item(5, 3).Value = 0.18
$ws.Cells.Item(5, 4).Value = 0.98

# --- Apply new number formats to the MAE (C) and R-squared (D) table columns ---
$ws.Columns.Item(3).NumberFormat = "0.000_ "
$ws.Columns.Item(4).NumberFormat = "0.000_);[Red]\(0.000\)"

# --- Make the R-squared header font white (matches poster styling) ---
$ws.Range("D1").Font.ThemeColor = 2

# --- Update view state: zoom in on the sheet and move selection to D1 ---
$excel.ActiveWindow.Zoom = 158
[void]$ws.Range("D1").Select()

Write-Host "done"
